$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: last_valuation_n -> last_valuation_b, exit_valuation_n -> exit_valuation_b
$ws.Range("H1").Value = "last_valuation_b"
$ws.Range("I1").Value = "exit_valuation_b"

# Convert columns H (last_valuation_n) and I (exit_valuation_n) from raw dollar
# amounts into billions of dollars for rows 2 through 49.
for ($r = 2; $r -le 49; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    if ($hCell.Value2 -ne $null) {
        $hCell.Value2 = $hCell.Value2 / 1000000000
    }
    $iCell = $ws.Cells.Item($r, 9)
    if ($iCell.Value2 -ne $null) {
        $iCell.Value2 = $iCell.Value2 / 1000000000
    }
}
